$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-7: ExecutionMode column flips from Yes to No ---
$ws.Cells.Item(3,4).Value = "No"
$ws.Cells.Item(4,4).Value = "No"
$ws.Cells.Item(5,4).Value = "No"
$ws.Cells.Item(6,4).Value = "No"
$ws.Cells.Item(7,4).Value = "No"

# --- New row 8: TestCase7 ---
$ws.Cells.Item(8,1).Value = "TestCase7"
$ws.Cells.Item(8,2).Value = "D:\SIN_ADE.pdf"
$ws.Cells.Item(8,3).Value = "D:\SIN_UW.pdf"
$ws.Cells.Item(8,4).Value = "No"
$ws.Cells.Item(8,1).WrapText = $true
$ws.Cells.Item(8,2).WrapText = $true
$ws.Cells.Item(8,3).WrapText = $true
$ws.Cells.Item(8,4).WrapText = $true

# --- New row 9: TestCase8 ---
$ws.Cells.Item(9,1).Value = "TestCase8"
$ws.Cells.Item(9,2).Value = "D:\Test.pdf"
$ws.Cells.Item(9,3).Value = "D:\Test.pdf"
$ws.Cells.Item(9,4).Value = "No"
$ws.Cells.Item(9,1).WrapText = $true
$ws.Cells.Item(9,4).WrapText = $true

# --- New row 10: TestCase9 ---
$ws.Cells.Item(10,1).Value = "TestCase9"
$ws.Cells.Item(10,2).Value = "D:\Test.pdf"
$ws.Cells.Item(10,3).Value = "D:\Test.pdf"
$ws.Cells.Item(10,4).Value = "No"
$ws.Cells.Item(10,1).WrapText = $true
$ws.Cells.Item(10,4).WrapText = $true

# --- New row 11: TestCase10 ---
$ws.Cells.Item(11,1).Value = "TestCase10"
$ws.Cells.Item(11,2).Value = "D:\Test.pdf"
$ws.Cells.Item(11,3).Value = "D:\Test.pdf"
$ws.Cells.Item(11,4).Value = "No"
$ws.Cells.Item(11,1).WrapText = $true
$ws.Cells.Item(11,4).WrapText = $true

# --- New row 12: TestCase11 ---
$ws.Cells.Item(12,1).Value = "TestCase11"
$ws.Cells.Item(12,2).Value = "D:\Test.pdf"
$ws.Cells.Item(12,3).Value = "D:\Test.pdf"
$ws.Cells.Item(12,4).Value = "No"
$ws.Cells.Item(12,1).WrapText = $true
$ws.Cells.Item(12,4).WrapText = $true

# --- Update row 2: InputFile1/InputFile2 now point to the new ABRCIR pdf ---
$ws.Cells.Item(2,2).Value = "D:\ABRCIR-20180326.pdf"
$ws.Cells.Item(2,3).Value = "D:\ABRCIR-20180326.pdf"

# --- Update selection to match the new active cell ---
$ws.Range("C2").Select() | Out-Null
